# Updated capital structure database
# Apply new computed metrics to the two data rows (row 2 and row 3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2,3) {
    $ws.Range("G$r").Value  = -0.01476510067114094
    $ws.Range("H$r").Value  = -0
    $ws.Range("I$r").Value  = -0.002042344702481116
    $ws.Range("J$r").Value  = -0.002042344702481116
    $ws.Range("K$r").Value  = -11.4
    $ws.Range("L$r").Value  = 3.825503355704698

    $ws.Range("W$r").Value  = -0.6745562130177516
    $ws.Range("X$r").Value  = 0.03890718064103192
    $ws.Range("Y$r").Value  = -0.7134633936587835
    $ws.Range("Z$r").Value  = -0.1405026189366338
    $ws.Range("AA$r").Value = 0.000286954779469957
    $ws.Range("AB$r").Value = 0.03383421247535916
    $ws.Range("AC$r").Value = -0.0335472576958892
    $ws.Range("AD$r").Value = 5.13
    $ws.Range("AE$r").Value = 0.3395690639330314
    $ws.Range("AF$r").Value = 5.469569063933031
    $ws.Range("AG$r").Value = 5.469569063933031
    $ws.Range("AH$r").Value = 0.3281170042822052
    $ws.Range("AI$r").Value = 0.5116734857336338
    $ws.Range("AJ$r").Value = 0.3281170042822052
    $ws.Range("AK$r").Value = 0.5116734857336338

    $ws.Range("AN$r").Value = 69.32432432432432
    $ws.Range("AP$r").Value = 73.91309545855448
}

# Row 2 specific values
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0
$ws.Range("S2").Value = 0
$ws.Range("T2").ClearContents()

# Row 3 specific values
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").ClearContents()
